$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter Salvar em pdf Salvar em docx"
$jupRange = $d.Content.Duplicate
$null = $jupRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false,
                                $false, $false, $true, 1, $false, "", 0)
$jupIndex = $jupRange.Paragraphs.First.Index

# Locate the paragraph that contains the trailing copyright notice
$copyRange = $d.Content.Duplicate
$null = $copyRange.Find.Execute("Creative Commons Attribution", $true, $false, $false,
                                 $false, $false, $true, 1, $false, "", 0)
$copyIndex = $copyRange.Paragraphs.First.Index

# Remove the blank paragraph right before "Ver no Jupiter..." through the end of the
# copyright paragraph (inclusive) - this drops the blank paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the "© 2020 ..." paragraph,
# while leaving the rest of the document (including the following blank paragraph and the
# final page-break paragraph) untouched.
$startParagraph = $d.Paragraphs.Item($jupIndex - 1)
$endParagraph = $d.Paragraphs.Item($copyIndex)

$deleteRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)
$deleteRange.Delete()
